# Auto-generated edit script applying numeric updates to Ragnarok_Profits sheets
# per the commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J2").Value = 1251080.4
$ws.Range("N2").Value = -1251306.4
$ws.Range("H2").Value = 1001404.3
$ws.Range("L2").Value = 1251080.4
$ws.Range("J6").Value = 240
$ws.Range("I6").Value = 190
$ws.Range("N6").Value = -944
$ws.Range("H6").Value = 210
$ws.Range("K6").Value = 570
$ws.Range("M6").Value = -458
$ws.Range("L6").Value = 720
$ws.Range("L12").Value = 965
$ws.Range("I12").Value = 359.33334
$ws.Range("K12").Value = 359.33334
$ws.Range("M12").Value = -189.33334
$ws.Range("N12").Value = -1305
$ws.Range("H12").Value = 662.1667
$ws.Range("J12").Value = 965
$ws.Range("H15").Value = 3909.3845
$ws.Range("M15").Value = -11559.1535
$ws.Range("I15").Value = 3909.3845
$ws.Range("K15").Value = 11728.1535
$ws.Range("J28").Value = 9992.875
$ws.Range("I28").Value = 390
$ws.Range("L28").Value = 9992.875
$ws.Range("K28").Value = 390
$ws.Range("H28").Value = 3881.9546
$ws.Range("M28").Value = 95
$ws.Range("N28").Value = -10962.875
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 30000
$ws.Range("H29").Value = 10000
$ws.Range("N29").Value = -30562
$ws.Range("H38").Value = 23709.666
$ws.Range("L38").Value = 19978.9995
$ws.Range("N38").Value = -20722.9995
$ws.Range("J38").Value = 6659.6665
$ws.Range("H86").Value = 6679.1113
$ws.Range("J86").Value = 10081.4
$ws.Range("L86").Value = 10081.4
$ws.Range("N86").Value = -12327.4
$ws.Range("H89").Value = 6679.1113
$ws.Range("N89").Value = -61639
$ws.Range("J89").Value = 10081.4
$ws.Range("L89").Value = 50407
$ws.Range("I107").Value = 694.6957
$ws.Range("M107").Value = 1225.3043
$ws.Range("K107").Value = 694.6957
$ws.Range("H107").Value = 927.7241
$ws.Range("L138").Value = 39576.3
$ws.Range("N138").Value = -49856.3
$ws.Range("H138").Value = 12146.923
$ws.Range("J138").Value = 13192.1
$ws.Range("H140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K141").Value = 300020190
$ws.Range("I141").Value = 100006730
$ws.Range("M141").Value = -300015010
$ws.Range("H141").Value = 62516460

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J2").Value = 4528.8184
$ws.Range("N2").Value = -4754.8184
$ws.Range("H2").Value = 5336.3335
$ws.Range("L2").Value = 4528.8184
$ws.Range("M32").Value = -3436.56
$ws.Range("K32").Value = 3723.56
$ws.Range("I32").Value = 3723.56
$ws.Range("H32").Value = 3723.56
$ws.Range("J43").Value = 38326.668
$ws.Range("H43").Value = 38674
$ws.Range("L43").Value = 38326.668
$ws.Range("N43").Value = -38952.668
$ws.Range("M110").Value = -6889.666999999999
$ws.Range("I110").Value = 8934.666999999999
$ws.Range("L110").Value = 5000
$ws.Range("J110").Value = 5000
$ws.Range("N110").Value = -9090
$ws.Range("K110").Value = 8934.666999999999
$ws.Range("H110").Value = 8372.571
$ws.Range("L116").Value = 4528.8184
$ws.Range("J116").Value = 4528.8184
$ws.Range("H116").Value = 5336.3335
$ws.Range("N116").Value = -9116.8184
$ws.Range("I132").Value = 4591.08
$ws.Range("H132").Value = 3576045.5
$ws.Range("K132").Value = 13773.24
$ws.Range("M132").Value = -11243.24

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L3").Value = 4528.8184
$ws.Range("J3").Value = 4528.8184
$ws.Range("H3").Value = 5336.3335
$ws.Range("N3").Value = -4756.8184
$ws.Range("M86").Value = -1151.1
$ws.Range("H86").Value = 4491.1875
$ws.Range("J86").Value = 8186.3335
$ws.Range("K86").Value = 2274.1
$ws.Range("L86").Value = 8186.3335
$ws.Range("I86").Value = 2274.1
$ws.Range("N86").Value = -10432.3335
$ws.Range("M89").Value = -5754.5
$ws.Range("I89").Value = 2274.1
$ws.Range("H89").Value = 4491.1875
$ws.Range("K89").Value = 11370.5
$ws.Range("N89").Value = -52163.6675
$ws.Range("J89").Value = 8186.3335
$ws.Range("L89").Value = 40931.6675
$ws.Range("H134").Value = 5002642
$ws.Range("M134").Value = -5555.117400000001
$ws.Range("K134").Value = 8090.117400000001
$ws.Range("I134").Value = 2696.7058

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 500
$ws.Range("M22").Value = -150
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M39").Value = -9608
$ws.Range("K39").Value = 9999
$ws.Range("I39").Value = 9999
$ws.Range("H39").Value = 24999
$ws.Range("H42").Value = 93000
$ws.Range("L42").Value = 93000
$ws.Range("J42").Value = 93000
$ws.Range("N42").Value = -94186
$ws.Range("M49").Value = -9817
$ws.Range("I49").Value = 9999
$ws.Range("K49").Value = 9999
$ws.Range("H49").Value = 24999
$ws.Range("L58").Value = 1999.6
$ws.Range("M58").Value = -2252.3547
$ws.Range("N58").Value = -2405.6
$ws.Range("I58").Value = 2455.3547
$ws.Range("J58").Value = 1999.6
$ws.Range("K58").Value = 2455.3547
$ws.Range("H58").Value = 2344.195
$ws.Range("M86").Value = -10115.467
$ws.Range("H86").Value = 11428.059
$ws.Range("J86").Value = 12850
$ws.Range("K86").Value = 11238.467
$ws.Range("L86").Value = 12850
$ws.Range("I86").Value = 11238.467
$ws.Range("N86").Value = -15096
$ws.Range("M89").Value = -50576.33500000001
$ws.Range("I89").Value = 11238.467
$ws.Range("H89").Value = 11428.059
$ws.Range("K89").Value = 56192.33500000001
$ws.Range("N89").Value = -75482
$ws.Range("J89").Value = 12850
$ws.Range("L89").Value = 64250
$ws.Range("H102").Value = 90332.664
$ws.Range("N102").Value = -119367
$ws.Range("J102").Value = 114499
$ws.Range("L102").Value = 114499
$ws.Range("I107").Value = 1121
$ws.Range("M107").Value = 799
$ws.Range("K107").Value = 1121
$ws.Range("H107").Value = 1525.6333
$ws.Range("H134").Value = 3742.4285
$ws.Range("M134").Value = -8692.2855
$ws.Range("K134").Value = 11227.2855
$ws.Range("I134").Value = 3742.4285
$ws.Range("H136").Value = 2344.195
$ws.Range("I136").Value = 2455.3547
$ws.Range("N136").Value = -11098.8
$ws.Range("K136").Value = 7366.0641
$ws.Range("M136").Value = -4816.0641
$ws.Range("J136").Value = 1999.6
$ws.Range("L136").Value = 5998.799999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 15321.255
$ws.Range("I56").Value = 15321.255
$ws.Range("M56").Value = -14791.255
$ws.Range("K56").Value = 15321.255
$ws.Range("H97").Value = 1348.3846
$ws.Range("N97").Value = -4559.5454
$ws.Range("L97").Value = 3567.5454
$ws.Range("J97").Value = 1189.1818
$ws.Range("J121").Value = 5253.067
$ws.Range("M121").Value = 110
$ws.Range("H121").Value = 4949.75
$ws.Range("K121").Value = 1200
$ws.Range("I121").Value = 400
$ws.Range("N121").Value = -18379.201
$ws.Range("L121").Value = 15759.201
$ws.Range("N131").Value = -33830.1432
$ws.Range("J131").Value = 7916.7144
$ws.Range("H131").Value = 6348.2856
$ws.Range("M131").Value = -9299.571
$ws.Range("I131").Value = 4779.857
$ws.Range("K131").Value = 14339.571
$ws.Range("L131").Value = 23750.1432

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2379.2307
$ws.Range("L80").Value = 3502.75
$ws.Range("M80").Value = -881.8888999999999
$ws.Range("K80").Value = 1879.8889
$ws.Range("J80").Value = 3502.75
$ws.Range("N80").Value = -5498.75
$ws.Range("I80").Value = 1879.8889
$ws.Range("H83").Value = 2379.2307
$ws.Range("L83").Value = 17513.75
$ws.Range("I83").Value = 1879.8889
$ws.Range("N83").Value = -27497.75
$ws.Range("M83").Value = -4407.4445
$ws.Range("J83").Value = 3502.75
$ws.Range("K83").Value = 9399.4445
$ws.Range("L113").Value = 9261009
$ws.Range("J113").Value = 9261009
$ws.Range("N113").Value = -9265349
$ws.Range("H113").Value = 9261009
$ws.Range("J122").Value = 2999
$ws.Range("H122").Value = 2284.3333
$ws.Range("L122").Value = 8997
$ws.Range("I122").Value = 2141.4
$ws.Range("M122").Value = -3974.200000000001
$ws.Range("N122").Value = -13897
$ws.Range("K122").Value = 6424.200000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J46").Value = 3993.3333
$ws.Range("H46").Value = 3494
$ws.Range("L46").Value = 3993.3333
$ws.Range("N46").Value = -4369.3333
$ws.Range("H122").Value = 3571.814
$ws.Range("I122").Value = 3280.6487
$ws.Range("M122").Value = -7391.946100000001
$ws.Range("K122").Value = 9841.946100000001
$ws.Range("I132").Value = 3332
$ws.Range("H132").Value = 5415.625
$ws.Range("K132").Value = 9996
$ws.Range("M132").Value = -7466

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K109").Value = 15000
$ws.Range("M109").Value = -13613
$ws.Range("I109").Value = 15000
$ws.Range("H109").Value = 57499.5
$ws.Range("H136").Value = 632499.3
$ws.Range("I136").Value = 8923.076999999999
$ws.Range("N136").Value = -10009089
$ws.Range("K136").Value = 26769.231
$ws.Range("M136").Value = -24219.231
$ws.Range("J136").Value = 3334663
$ws.Range("L136").Value = 10003989
$ws.Range("N137").Value = -104949.5
$ws.Range("J137").Value = 94749.5
$ws.Range("L137").Value = 94749.5
$ws.Range("H137").Value = 94749.5

# ---- Cell removals (value cleared entirely, not set to 0) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N140").ClearContents()
